# Fix avg/std experiment-time related columns (G, H, M, N) for rows 2-13
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 6.672440720000001
$ws.Range("H2").Value = 381.69752676
$ws.Range("M2").Value = 0.494969513974153
$ws.Range("N2").Value = 50.73986853472349
$ws.Range("G3").Value = 6.899435
$ws.Range("H3").Value = 627.70538312
$ws.Range("M3").Value = 0.6730628938948248
$ws.Range("N3").Value = 118.6616133725753
$ws.Range("G4").Value = 2.58612274
$ws.Range("H4").Value = 80.25798019999999
$ws.Range("M4").Value = 0.3873316922172573
$ws.Range("N4").Value = 22.46647487508017
$ws.Range("G5").Value = 2.77142266
$ws.Range("H5").Value = 136.06397188
$ws.Range("M5").Value = 0.3576325808935163
$ws.Range("N5").Value = 32.5564593688892
$ws.Range("G6").Value = 0.86578674
$ws.Range("H6").Value = 14.07402488
$ws.Range("M6").Value = 0.2366681157342128
$ws.Range("N6").Value = 6.41590886572101
$ws.Range("G7").Value = 0.9817642000000001
$ws.Range("H7").Value = 25.13027774
$ws.Range("M7").Value = 0.1766459389475266
$ws.Range("N7").Value = 7.867535293885537
$ws.Range("G8").Value = 0.40571482
$ws.Range("H8").Value = 4.064102399999999
$ws.Range("M8").Value = 0.1366649482916183
$ws.Range("N8").Value = 2.15126537033222
$ws.Range("G9").Value = 0.485242
$ws.Range("H9").Value = 8.82604478
$ws.Range("M9").Value = 0.1230362935864006
$ws.Range("N9").Value = 4.355758105270057
$ws.Range("G10").Value = 0.22744752
$ws.Range("H10").Value = 1.74933596
$ws.Range("M10").Value = 0.09301120797337778
$ws.Range("N10").Value = 1.044050193472449
$ws.Range("G11").Value = 0.25713032
$ws.Range("H11").Value = 3.54976268
$ws.Range("M11").Value = 0.07773953739795017
$ws.Range("N11").Value = 1.991305742254033
$ws.Range("G12").Value = 0.13371546
$ws.Range("H12").Value = 0.8568562200000001
$ws.Range("M12").Value = 0.06095304612940445
$ws.Range("N12").Value = 0.5810678937613333
$ws.Range("G13").Value = 0.15202994
$ws.Range("H13").Value = 1.7709933
$ws.Range("M13").Value = 0.05048040274322935
$ws.Range("N13").Value = 1.07030705865079
